# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to currentAveragePrice / Leve price / profit columns
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4819.6
$ws.Range("J29").Value = 5998.75
$ws.Range("L29").Value = 17996.25
$ws.Range("N29").Value = -18558.25
$ws.Range("H92").Value = 4172.3335
$ws.Range("I92").Value = 7012
$ws.Range("J92").Value = 1332.6666
$ws.Range("K92").Value = 7012
$ws.Range("L92").Value = 1332.6666
$ws.Range("M92").Value = -5764
$ws.Range("N92").Value = -3828.6666
$ws.Range("H97").Value = 2391
$ws.Range("J97").Value = 2391
$ws.Range("L97").Value = 7173
$ws.Range("N97").Value = -8165
$ws.Range("H99").Value = 229.75
$ws.Range("I99").Value = 233.33333
$ws.Range("J99").Value = 219
$ws.Range("K99").Value = 699.99999
$ws.Range("L99").Value = 657
$ws.Range("M99").Value = 798.00001
$ws.Range("N99").Value = -3653
$ws.Range("H106").Value = 12229.7
$ws.Range("I106").Value = 14162.25
$ws.Range("J106").Value = 4499.5
$ws.Range("K106").Value = 14162.25
$ws.Range("L106").Value = 4499.5
$ws.Range("M106").Value = -13531.25
$ws.Range("N106").Value = -5761.5
$ws.Range("H138").Value = 4911.16
$ws.Range("I138").Value = 6499
$ws.Range("J138").Value = 4773.087
$ws.Range("K138").Value = 19497
$ws.Range("L138").Value = 14319.261
$ws.Range("M138").Value = -14357
$ws.Range("N138").Value = -24599.261

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2646826.8
$ws.Range("I2").Value = 3473172.5
$ws.Range("K2").Value = 3473172.5
$ws.Range("M2").Value = -3473059.5
$ws.Range("H45").Value = 4798078.5
$ws.Range("I45").Value = 6255959.5
$ws.Range("K45").Value = 6255959.5
$ws.Range("M45").Value = -6255582.5
$ws.Range("H61").Value = 8558.235000000001
$ws.Range("I61").Value = 9680.786
$ws.Range("J61").Value = 3319.6667
$ws.Range("K61").Value = 9680.786
$ws.Range("L61").Value = 3319.6667
$ws.Range("M61").Value = -9468.786
$ws.Range("N61").Value = -3743.6667
$ws.Range("H110").Value = 2138358.8
$ws.Range("I110").Value = 6947297
$ws.Range("K110").Value = 6947297
$ws.Range("M110").Value = -6945252
$ws.Range("H116").Value = 2646826.8
$ws.Range("I116").Value = 3473172.5
$ws.Range("K116").Value = 3473172.5
$ws.Range("M116").Value = -3470878.5
$ws.Range("H132").Value = 6282.242
$ws.Range("I132").Value = 6430.4585
$ws.Range("J132").Value = 5887
$ws.Range("K132").Value = 19291.3755
$ws.Range("L132").Value = 17661
$ws.Range("M132").Value = -16761.3755
$ws.Range("N132").Value = -22721
$ws.Range("H136").Value = 8558.235000000001
$ws.Range("I136").Value = 9680.786
$ws.Range("J136").Value = 3319.6667
$ws.Range("K136").Value = 29042.358
$ws.Range("L136").Value = 9959.000100000001
$ws.Range("M136").Value = -26492.358
$ws.Range("N136").Value = -15059.0001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2646826.8
$ws.Range("I3").Value = 3473172.5
$ws.Range("K3").Value = 3473172.5
$ws.Range("M3").Value = -3473058.5
$ws.Range("H24").Value = 2016
$ws.Range("I24").Value = 2016
$ws.Range("K24").Value = 2016
$ws.Range("M24").Value = -1781
$ws.Range("H61").Value = 88553.5
$ws.Range("J61").Value = 88553.5
$ws.Range("L61").Value = 88553.5
$ws.Range("N61").Value = -89179.5
$ws.Range("H76").Value = 21314
$ws.Range("J76").Value = 21314
$ws.Range("L76").Value = 21314
$ws.Range("N76").Value = -21944
$ws.Range("H79").Value = 21314
$ws.Range("J79").Value = 21314
$ws.Range("L79").Value = 21314
$ws.Range("N79").Value = -23498
$ws.Range("H94").Value = 2020961.2
$ws.Range("I94").Value = 2674351
$ws.Range("J94").Value = 1393.091
$ws.Range("K94").Value = 2674351
$ws.Range("L94").Value = 1393.091
$ws.Range("M94").Value = -2673900
$ws.Range("N94").Value = -2295.091
$ws.Range("H99").Value = 3487189.2
$ws.Range("I99").Value = 5497230
$ws.Range("J99").Value = 3118.0667
$ws.Range("K99").Value = 5497230
$ws.Range("L99").Value = 3118.0667
$ws.Range("M99").Value = -5495732
$ws.Range("N99").Value = -6114.066699999999
$ws.Range("H134").Value = 11661.37
$ws.Range("I134").Value = 11266.5
$ws.Range("J134").Value = 13398.8
$ws.Range("K134").Value = 33799.5
$ws.Range("L134").Value = 40196.39999999999
$ws.Range("M134").Value = -31264.5
$ws.Range("N134").Value = -45266.39999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 596.1429000000001
$ws.Range("I22").Value = 528.6667
$ws.Range("K22").Value = 528.6667
$ws.Range("M22").Value = -178.6667
$ws.Range("H122").Value = 2149.4285
$ws.Range("I122").Value = 1912.9333
$ws.Range("J122").Value = 2740.6667
$ws.Range("K122").Value = 5738.7999
$ws.Range("L122").Value = 8222.000100000001
$ws.Range("M122").Value = -3288.7999
$ws.Range("N122").Value = -13122.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 98.875
$ws.Range("I10").Value = 111.57143
$ws.Range("K10").Value = 334.71429
$ws.Range("M10").Value = -195.71429
$ws.Range("H80").Value = 2747
$ws.Range("I80").Value = 2374.25
$ws.Range("J80").Value = 3492.5
$ws.Range("K80").Value = 7122.75
$ws.Range("L80").Value = 10477.5
$ws.Range("M80").Value = -6186.75
$ws.Range("N80").Value = -12349.5
$ws.Range("H83").Value = 2747
$ws.Range("I83").Value = 2374.25
$ws.Range("J83").Value = 3492.5
$ws.Range("K83").Value = 21368.25
$ws.Range("L83").Value = 31432.5
$ws.Range("M83").Value = -16688.25
$ws.Range("N83").Value = -40792.5
$ws.Range("H132").Value = 1705.25
$ws.Range("J132").Value = 2248
$ws.Range("L132").Value = 20232
$ws.Range("N132").Value = -25292

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 39983.668
$ws.Range("J39").Value = 39983.668
$ws.Range("L39").Value = 39983.668
$ws.Range("N39").Value = -41047.668
$ws.Range("H70").Value = 20005452
$ws.Range("I70").Value = 28576074
$ws.Range("J70").Value = 7333
$ws.Range("K70").Value = 28576074
$ws.Range("L70").Value = 7333
$ws.Range("M70").Value = -28575804
$ws.Range("N70").Value = -7873
$ws.Range("H73").Value = 20005452
$ws.Range("I73").Value = 28576074
$ws.Range("J73").Value = 7333
$ws.Range("K73").Value = 28576074
$ws.Range("L73").Value = 7333
$ws.Range("M73").Value = -28575138
$ws.Range("N73").Value = -9205
$ws.Range("H97").Value = 1135207.9
$ws.Range("I97").Value = 2382564.2
$ws.Range("J97").Value = 1247.5454
$ws.Range("K97").Value = 2382564.2
$ws.Range("L97").Value = 1247.5454
$ws.Range("M97").Value = -2382068.2
$ws.Range("N97").Value = -2239.5454
$ws.Range("H101").Value = 23300
$ws.Range("J101").Value = 23300
$ws.Range("L101").Value = 23300
$ws.Range("N101").Value = -29790
$ws.Range("H102").Value = 5036991
$ws.Range("I102").Value = 7409177.5
$ws.Range("J102").Value = 2071757.4
$ws.Range("K102").Value = 7409177.5
$ws.Range("L102").Value = 2071757.4
$ws.Range("M102").Value = -7407555.5
$ws.Range("N102").Value = -2075001.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1036.4445
$ws.Range("I16").Value = 582.5714
$ws.Range("K16").Value = 582.5714
$ws.Range("M16").Value = -412.5714
$ws.Range("H40").Value = 10452.909
$ws.Range("I40").Value = 9998.299999999999
$ws.Range("K40").Value = 9998.299999999999
$ws.Range("M40").Value = -9862.299999999999
$ws.Range("H75").Value = 40076
$ws.Range("I75").Value = 20157
$ws.Range("J75").Value = 59995
$ws.Range("K75").Value = 20157
$ws.Range("L75").Value = 59995
$ws.Range("M75").Value = -19221
$ws.Range("N75").Value = -61867
$ws.Range("H78").Value = 40076
$ws.Range("I78").Value = 20157
$ws.Range("J78").Value = 59995
$ws.Range("K78").Value = 60471
$ws.Range("L78").Value = 179985
$ws.Range("M78").Value = -55791
$ws.Range("N78").Value = -189345
$ws.Range("H93").Value = 12346954
$ws.Range("I93").Value = 15152671
$ws.Range("J93").Value = 1800
$ws.Range("K93").Value = 15152671
$ws.Range("L93").Value = 1800
$ws.Range("M93").Value = -15151423
$ws.Range("N93").Value = -4296
$ws.Range("H136").Value = 37558.2
$ws.Range("I136").Value = 53112.25
$ws.Range("J136").Value = 6450.1
$ws.Range("K136").Value = 159336.75
$ws.Range("L136").Value = 19350.3
$ws.Range("M136").Value = -156786.75
$ws.Range("N136").Value = -24450.3

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H100").Value = 906.1111
$ws.Range("I100").Value = 360.83334
$ws.Range("K100").Value = 721.66668
$ws.Range("M100").Value = -180.66668
$ws.Range("H107").Value = 40002196
$ws.Range("I107").Value = 47621532
$ws.Range("K107").Value = 142864596
$ws.Range("M107").Value = -142862676
$ws.Range("H122").Value = 1697.091
$ws.Range("I122").Value = 1499.2222
$ws.Range("K122").Value = 4497.6666
$ws.Range("M122").Value = -2047.6666
$ws.Range("H126").Value = 2897.25
$ws.Range("I126").Value = 2776.9
$ws.Range("J126").Value = 3499
$ws.Range("K126").Value = 8330.700000000001
$ws.Range("L126").Value = 10497
$ws.Range("M126").Value = -5860.700000000001
$ws.Range("N126").Value = -15437
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()

